$d = $word.ActiveDocument
Write-Output $d.Content.Text
